$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("A4").Value = "ECs"
$ws.Range("A5").Value = "ECs"
$ws.Range("A6").Value = "FAPs"
$ws.Range("A7").Value = "FAPs"
$ws.Range("A8").Value = "FAPs"
$ws.Range("A9").Value = "FAPs"
$ws.Range("A10").Value = "sCs"
$ws.Range("A11").Value = "sCs"
$ws.Range("A12").Value = "sCs"
$ws.Range("A13").Value = "sCs"
$ws.Range("B2").Value = "Lgi2"
$ws.Range("B3").Value = "Lgi2"
$ws.Range("B4").Value = "Lgi2"
$ws.Range("B5").Value = "Lgi2"
$ws.Range("B6").Value = "Lgi2"
$ws.Range("B7").Value = "Lgi2"
$ws.Range("B8").Value = "Lgi2"
$ws.Range("B9").Value = "Lgi2"
$ws.Range("B10").Value = "Lgi2"
$ws.Range("B11").Value = "Lgi2"
$ws.Range("B12").Value = "Lgi2"
$ws.Range("B13").Value = "Lgi2"
$ws.Range("C2").Value = "Adam23"
$ws.Range("C3").Value = "Adam23"
$ws.Range("C4").Value = "Adam23"
$ws.Range("C5").Value = "Adam23"
$ws.Range("C6").Value = "Adam23"
$ws.Range("C7").Value = "Adam23"
$ws.Range("C8").Value = "Adam23"
$ws.Range("C9").Value = "Adam23"
$ws.Range("C10").Value = "Adam23"
$ws.Range("C11").Value = "Adam23"
$ws.Range("C12").Value = "Adam23"
$ws.Range("C13").Value = "Adam23"
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "M2"
$ws.Range("D5").Value = "sCs"
$ws.Range("D6").Value = "ECs"
$ws.Range("D7").Value = "FAPs"
$ws.Range("D8").Value = "M2"
$ws.Range("D9").Value = "sCs"
$ws.Range("D10").Value = "ECs"
$ws.Range("D11").Value = "FAPs"
$ws.Range("D12").Value = "M2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02053466666666667
$ws.Range("H2").Value = 0.061604
$ws.Range("I2").Value = 0.001930134604048101
$ws.Range("J2").Value = 0.001930134604048101
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1448543333333333
$ws.Range("N2").Value = 0.434563
$ws.Range("O2").Value = 0.007514552731137001
$ws.Range("P2").Value = 0.007514552731137
$ws.Range("Q2").Value = 0.002974535450222222
$ws.Range("R2").Value = 0.026770819052
$ws.Range("S2").Value = 0.00001450409826031169
$ws.Range("T2").Value = 0.00001450409826031169
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02053466666666667
$ws.Range("H3").Value = 0.061604
$ws.Range("I3").Value = 0.001930134604048101
$ws.Range("J3").Value = 0.001930134604048101
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 14.494489
$ws.Range("N3").Value = 43.483467
$ws.Range("O3").Value = 0.7519250504625467
$ws.Range("P3").Value = 0.7519250504625465
$ws.Range("Q3").Value = 0.2976395001186667
$ws.Range("R3").Value = 2.678755501068
$ws.Range("S3").Value = 0.001451316559548376
$ws.Range("T3").Value = 0.001451316559548375
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02053466666666667
$ws.Range("H4").Value = 0.061604
$ws.Range("I4").Value = 0.001930134604048101
$ws.Range("J4").Value = 0.001930134604048101
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.091182
$ws.Range("N4").Value = 0.273546
$ws.Range("O4").Value = 0.004730213666123443
$ws.Range("P4").Value = 0.004730213666123443
$ws.Range("Q4").Value = 0.001872391976
$ws.Range("R4").Value = 0.016851527784
$ws.Range("S4").Value = 0.000009129949081526087
$ws.Range("T4").Value = 0.000009129949081526087
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.02053466666666667
$ws.Range("H5").Value = 0.061604
$ws.Range("I5").Value = 0.001930134604048101
$ws.Range("J5").Value = 0.001930134604048101
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.545982333333334
$ws.Range("N5").Value = 13.637947
$ws.Range("O5").Value = 0.235830183140193
$ws.Range("P5").Value = 0.2358301831401929
$ws.Range("Q5").Value = 0.09335023188755556
$ws.Range("R5").Value = 0.840152086988
$ws.Range("S5").Value = 0.0004551839971578874
$ws.Range("T5").Value = 0.0004551839971578873
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 10.27464766666667
$ws.Range("H6").Value = 30.823943
$ws.Range("I6").Value = 0.9657548051669733
$ws.Range("J6").Value = 0.9657548051669733
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1448543333333333
$ws.Range("N6").Value = 0.434563
$ws.Range("O6").Value = 0.007514552731137001
$ws.Range("P6").Value = 0.007514552731137
$ws.Range("Q6").Value = 1.488327237989889
$ws.Range("R6").Value = 13.394945141909
$ws.Range("S6").Value = 0.007257215408776161
$ws.Range("T6").Value = 0.007257215408776161
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.27464766666667
$ws.Range("H7").Value = 30.823943
$ws.Range("I7").Value = 0.9657548051669733
$ws.Range("J7").Value = 0.9657548051669733
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 14.494489
$ws.Range("N7").Value = 43.483467
$ws.Range("O7").Value = 0.7519250504625467
$ws.Range("P7").Value = 0.7519250504625465
$ws.Range("Q7").Value = 148.9257675833757
$ws.Range("R7").Value = 1340.331908250381
$ws.Range("S7").Value = 0.7261752306096233
$ws.Range("T7").Value = 0.7261752306096232
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.27464766666667
$ws.Range("H8").Value = 30.823943
$ws.Range("I8").Value = 0.9657548051669733
$ws.Range("J8").Value = 0.9657548051669733
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.091182
$ws.Range("N8").Value = 0.273546
$ws.Range("O8").Value = 0.004730213666123443
$ws.Range("P8").Value = 0.004730213666123443
$ws.Range("Q8").Value = 0.936862923542
$ws.Range("R8").Value = 8.431766311878
$ws.Range("S8").Value = 0.0045682265775252
$ws.Range("T8").Value = 0.0045682265775252
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.27464766666667
$ws.Range("H9").Value = 30.823943
$ws.Range("I9").Value = 0.9657548051669733
$ws.Range("J9").Value = 0.9657548051669733
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.545982333333334
$ws.Range("N9").Value = 13.637947
$ws.Range("O9").Value = 0.235830183140193
$ws.Range("P9").Value = 0.2358301831401929
$ws.Range("Q9").Value = 46.70836677389123
$ws.Range("R9").Value = 420.375300965021
$ws.Range("S9").Value = 0.2277541325710487
$ws.Range("T9").Value = 0.2277541325710487
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.3437993333333333
$ws.Range("H10").Value = 1.031398
$ws.Range("I10").Value = 0.03231506022897868
$ws.Range("J10").Value = 0.03231506022897868
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.1448543333333333
$ws.Range("N10").Value = 0.434563
$ws.Range("O10").Value = 0.007514552731137001
$ws.Range("P10").Value = 0.007514552731137
$ws.Range("Q10").Value = 0.04980082323044444
$ws.Range("R10").Value = 0.448207409074
$ws.Range("S10").Value = 0.0002428332241005284
$ws.Range("T10").Value = 0.0002428332241005284
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.3437993333333333
$ws.Range("H11").Value = 1.031398
$ws.Range("I11").Value = 0.03231506022897868
$ws.Range("J11").Value = 0.03231506022897868
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 14.494489
$ws.Range("N11").Value = 43.483467
$ws.Range("O11").Value = 0.7519250504625467
$ws.Range("P11").Value = 0.7519250504625465
$ws.Range("Q11").Value = 4.983195655207333
$ws.Range("R11").Value = 44.848760896866
$ws.Range("S11").Value = 0.02429850329337503
$ws.Range("T11").Value = 0.02429850329337502
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.3437993333333333
$ws.Range("H12").Value = 1.031398
$ws.Range("I12").Value = 0.03231506022897868
$ws.Range("J12").Value = 0.03231506022897868
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.091182
$ws.Range("N12").Value = 0.273546
$ws.Range("O12").Value = 0.004730213666123443
$ws.Range("P12").Value = 0.004730213666123443
$ws.Range("Q12").Value = 0.03134831081199999
$ws.Range("R12").Value = 0.282134797308
$ws.Range("S12").Value = 0.0001528571395167171
$ws.Range("T12").Value = 0.0001528571395167171
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.3437993333333333
$ws.Range("H13").Value = 1.031398
$ws.Range("I13").Value = 0.03231506022897868
$ws.Range("J13").Value = 0.03231506022897868
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 4.545982333333334
$ws.Range("N13").Value = 13.637947
$ws.Range("O13").Value = 0.235830183140193
$ws.Range("P13").Value = 0.2358301831401929
$ws.Range("Q13").Value = 1.562905695545111
$ws.Range("R13").Value = 14.066151259906
$ws.Range("S13").Value = 0.007620866571986409
$ws.Range("T13").Value = 0.007620866571986406
